# Update column F (dSF) values for rows 8, 11, 15, 16, 23
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F8").Value = -9
$ws.Range("F11").Value = 3
$ws.Range("F15").Value = -1
$ws.Range("F16").Value = -6
$ws.Range("F23").Value = 1
